# Applies the cryptos-list price/volume refresh described by the commit diff.
# D-column ("Price") values are written as literal text: several look like plain
# numbers ("232.36") or have multiple dots ("43.305.61") and must stay text (as the
# original inlineStr cells were), so each D write is bracketed with a temporary
# NumberFormat = "@" (Text) that is reverted via Style = "Normal" afterwards - this
# keeps the cell text-typed without leaving a new/changed number format on it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.305.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.27%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.363.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.91%  "

# Row 4
$ws.Range("E4").Value = "  -0.88%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.36"
$ws.Range("D5").Style = "Normal"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.650"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.38%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "68.12"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.91%  "

# Row 8
$ws.Range("E8").Value = "  -0.05%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.459"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.76%  "

# Row 10
$ws.Range("E10").Value = "  -2.85%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.86"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.18%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.51"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.53%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.716.92"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.87%  "

# Row 14
$ws.Range("E14").Value = "  -1.46%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.43%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.10%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.842"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.04%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.368.61"
$ws.Range("D18").Style = "Normal"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.374.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.95%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0978"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.37%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.03%  "

# Row 22
$ws.Range("E22").Value = "  +3.02%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "248.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.89%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +16.61%  "

# Row 25
$ws.Range("E25").Value = "  -0.03%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.50%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.02%  "

# Row 28
$ws.Range("E28").Value = "  -1.59%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.86%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "174.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.94%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.54"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.78%  "

# Row 32
$ws.Range("E32").Value = "  -6.86%  "

# Row 33
$ws.Range("E33").Value = "  -0.10%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.98"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.63%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0691"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.10%  "

# Row 36
$ws.Range("E36").Value = "  +3.17%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.77%  "

# Row 38
$ws.Range("E38").Value = "  +1.06%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.85%  "

# Row 40
$ws.Range("E40").Value = "  -1.88%  "

# Row 41
$ws.Range("E41").Value = "  -0.02%  "

# Row 42
$ws.Range("E42").Value = "  +8.53%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "18.16"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.67%  "

# Row 44
$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.50"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.63%  "

# Row 45
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.81%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "98.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.33%  "

# Row 47
$ws.Range("E47").Value = "  +1.50%  "

# Row 48
$ws.Range("E48").Value = "  +0.24%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.447.19"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.79%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.589.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.07%  "

# Row 51
$ws.Range("E51").Value = "  -2.87%  "
